$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# The universal separator is now a comma instead of a plus sign.
$ws.Range("B2").Value = "Wall, Floor"

# Update the active selection as recorded in the sheet view.
$ws.Range("B8").Select()
